$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for rows 2-56 (A:POS, B:nSubjects, C:nSamples, D:gene, E:mutation)
$data = @(
    @(241, 30, 56, 'transgenic', ' '),
    @(3037, 30, 56, 'orf1ab', 'silent'),
    @(8782, 30, 56, 'orf1ab', 'silent'),
    @(14408, 30, 56, 'orf1ab', 'P314L'),
    @(18060, 30, 56, 'orf1ab', 'silent'),
    @(23403, 30, 56, 'S', 'D614G'),
    @(28144, 30, 56, 'orf8', 'S84L'),
    @(25563, 28, 51, 'orf3a', 'Q57H'),
    @(1059, 23, 40, 'orf1ab', 'T265I'),
    @(16260, 7, 11, 'orf1ab', 'silent'),
    @(28821, 7, 11, 'N', 'S183Y'),
    @(18255, 5, 10, 'orf1ab', 'M1596I'),
    @(18877, 4, 10, 'orf1ab', 'silent'),
    @(11916, 4, 9, 'orf1ab', 'S3884L'),
    @(21986, 5, 9, 'S', 'G142S'),
    @(18998, 3, 7, 'orf1ab', 'A1844V'),
    @(29540, 3, 7, 'transgenic', ' '),
    @(1918, 2, 6, 'orf1ab', 'silent'),
    @(23230, 2, 6, 'S', 'silent'),
    @(7386, 1, 5, 'orf1ab', 'M2374T'),
    @(11083, 3, 5, 'orf1ab', 'L3606F'),
    @(23994, 1, 5, 'S', 'K811T'),
    @(1500, 2, 4, 'orf1ab', 'S412Y'),
    @(3090, 2, 4, 'orf1ab', 'E942G'),
    @(4113, 2, 4, 'orf1ab', 'A1283V'),
    @(13501, 2, 4, 'orf1ab', 'P12S'),
    @(20371, 2, 4, 'orf1ab', 'silent'),
    @(24197, 2, 4, 'S', 'A879S'),
    @(28109, 2, 4, 'orf8', 'Q72H'),
    @(29344, 2, 4, 'N', 'silent'),
    @(6388, 1, 3, 'orf1ab', 'silent'),
    @(15326, 2, 3, 'orf1ab', 'M620T'),
    @(28609, 1, 3, 'N', 'silent'),
    @(29247, 1, 3, 'N', 'T325I'),
    @(157, 1, 2, 'transgenic', ' '),
    @(2219, 1, 2, 'orf1ab', 'I652V'),
    @(2485, 1, 2, 'orf1ab', 'silent'),
    @(4002, 1, 2, 'orf1ab', 'T1246I'),
    @(5570, 2, 2, 'orf1ab', 'M1769V'),
    @(12513, 1, 2, 'orf1ab', 'T4083M'),
    @(13975, 1, 2, 'orf1ab', 'G170C'),
    @(15857, 1, 2, 'orf1ab', 'T797I'),
    @(17010, 1, 2, 'orf1ab', 'silent'),
    @(18512, 2, 2, 'orf1ab', 'P1682L'),
    @(18650, 1, 2, 'orf1ab', 'E1728V'),
    @(22317, 1, 2, 'S', 'G252D'),
    @(28068, 1, 2, 'orf8', 'E59*'),
    @(29614, 1, 2, 'orf10', 'silent'),
    @(625, 1, 1, 'orf1ab', 'K120N'),
    @(1403, 1, 1, 'orf1ab', 'P380S'),
    @(2262, 1, 1, 'orf1ab', 'T666I'),
    @(2334, 1, 1, 'orf1ab', 'A690V'),
    @(2416, 1, 1, 'orf1ab', 'silent'),
    @(2788, 1, 1, 'orf1ab', 'silent'),
    @(3312, 1, 1, 'orf1ab', 'L1016P')
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}
